$d = $word.ActiveDocument

# Helper: insert a WordprocessingML fragment (wrapped as a mini OOXML package)
# immediately before the document's trailing paragraph, using the always
# up-to-date "last paragraph" range so repeated calls stack in order.
function InsertBodyXml($innerXml) {
    $p = $d.Paragraphs.Item($d.Paragraphs.Count)
    $r = $p.Range
    $r.Collapse(1)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           $innerXml +
           '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# 1. Register the "Hyperlink" character style the way Word does the first
#    time a hyperlink is inserted into a document that doesn't have it yet.
$hlStyle = $d.Styles.Add("Collegamento ipertestuale", 2)
$hlStyle.NameLocal = "Hyperlink"
$hlStyle.BaseStyle = $d.Styles.Item("Carpredefinitoparagrafo")
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.Font.Color = 16711680
$hlStyle.Font.Underline = 1

# 2. New bullet: "(ask about caching system)" - continues the existing
#    list (numId 4) right after "...calls ws"
InsertBodyXml '<w:body><w:p><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>(ask about caching system)</w:t></w:r></w:p></w:body>'

# 3. Hyperlink paragraph: http://goo.gl/x2Wbw1
InsertBodyXml '<w:body><w:p></w:p></w:body>'
$hostPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$hyperlink = $d.Hyperlinks.Add($hostPara.Range, "http://goo.gl/x2Wbw1", [Type]::Missing, [Type]::Missing, "http://goo.gl/x2Wbw1")
$hyperlink.Range.Style = $hlStyle

# 4. "response listener" paragraph
InsertBodyXml '<w:body><w:p><w:r><w:t>response listener</w:t></w:r></w:p></w:body>'

# 5. "stringrequest instead of jsonrequest since not all servers reply using json"
InsertBodyXml ('<w:body><w:p>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>stringrequest</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> instead of </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>jsonrequest</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> since not all servers reply using </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>json</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '</w:p></w:body>')
